$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# File-name cells (column D = Neo4jData file, column E = WebData file)
# TC01 / DNBSEQ-G400 -> TC12 / NextSeq 500
# ---------------------------------------------------------------------------
$neo4jName = 'TC12_CDS_Filter_InstrumentModel-NextSeq 500_Neo4jData.xlsx'
$webName   = 'TC12_CDS_Filter_InstrumentModel-NextSeq 500_WebData.xlsx'

$ws.Range("D2").Value = $neo4jName
$ws.Range("D3").Value = $neo4jName
$ws.Range("D4").Value = $neo4jName

$ws.Range("E2").Value = $webName
$ws.Range("E3").Value = $webName
$ws.Range("E4").Value = $webName

# ---------------------------------------------------------------------------
# Query text cells - swap the filtered instrument model from DNBSEQ-G400
# to "NextSeq 500" while preserving the Cypher text verbatim.
# ---------------------------------------------------------------------------

# B2 - ParticipantsTab query (Participant ID / Study Name / Accession / Gender / Samples)
$participantQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['NextSeq 500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@

# B3 - SamplesTab query (Sample ID / Participant ID / Study Name / Accession / Tumor / Analyte Type)
$sampleQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['NextSeq 500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

# B4 - FilesTab query (File Name / Study Name / Accession / Participant ID / Sample ID / File Type)
$fileQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['NextSeq 500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

# C2/C3/C4 - StatQuery (Studies / Participants / Samples / Files counts)
$statQuery = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['NextSeq 500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@

$ws.Range("B2").Value = $participantQuery
$ws.Range("B3").Value = $sampleQuery
$ws.Range("B4").Value = $fileQuery

$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# ---------------------------------------------------------------------------
# Column D widened (auto-fit grew slightly because of the new label text)
# and the active selection moved from B4 to D4.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 84.59244791666667

[void]$ws.Range("D4").Select()
